# Capitalize the first letter of the header cells in rows 14, 17, and 23
# (columns A/B/C), matching the source-language wording update:
#   "жынысы боюнча" / "по полу" / "by sex"                       -> capitalized
#   "энесинин билими " / "образование матери " / "education..." -> capitalized
#   "квинтиль по индексу благосостояния" / "wealth quintile"     -> capitalized
# (row 23's column A was already capitalized, so it is left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Жынысы боюнча"
$ws.Range("B14").Value = "По полу"
$ws.Range("B17").Value = "Образование матери "
$ws.Range("A17").Value = "Энесинин билими "
$ws.Range("B23").Value = "Квинтиль по индексу благосостояния"
$ws.Range("C14").Value = "By sex"
$ws.Range("C17").Value = "Education of mother"
$ws.Range("C23").Value = "Wealth quintile"

# Clear the active-cell selection that had been persisted in the sheet view
# (so the sheet no longer reports a stale A23 selection).
$ws.Range("A1").Select()
